$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table originally listed 36 parameter combinations covering two models
# (GWAS, RR-BLUP) crossed with two normalizations (Minmax applied twice: once
# for the 0.05-0.25 cut-offs, once for the 0.5-1 cut-offs). For the journal
# submission, the RR-BLUP rows are dropped entirely, and the surviving GWAS
# block that used to repeat "Minmax" for cut-offs 0.5-1 is now re-labeled as
# the "Z-score" normalization condition.

# Step 1: preserve the bottom border formatting that lived on the last row
# (row 37) by copying it onto row 27 - the row that will become the new last
# row (19) once the intervening RR-BLUP rows are removed.
$srcFormat = $ws.Range("A37:D37")
$srcFormat.Copy()
$dstFormat = $ws.Range("A27:D27")
$dstFormat.PasteSpecial(-4122)  # xlPasteFormats

# Step 2: remove the second RR-BLUP block (combinations 27-36, rows 28-37).
$ws.Rows("28:37").Delete()

# Step 3: remove the first RR-BLUP block (combinations 9-16, rows 10-17).
$ws.Rows("10:17").Delete()

# Step 4: the remaining rows 10-19 (formerly combinations 17-26, GWAS/Minmax
# 0.5-1) are renumbered to combinations 9-18 and relabeled with the new
# "Z-score" normalization.
for ($r = 10; $r -le 19; $r++) {
  $ws.Cells.Item($r, 1).Value = $r - 1
  $ws.Cells.Item($r, 3).Value = "Z-score"
}

# Step 5: restore the cursor/selection position recorded in the saved file.
$ws.Range("H15").Select() | Out-Null
